# Applies the "spreadsheet and tooltiptext improvement" edit:
#  - fieldnames: rewrite header/row2 labels, clear rows 3-6 (old
#    Drive/Gapminder/Google Datasearch/Kaggle/OurWorldInData leftovers),
#    drop their hyperlinks (keep the Expasy one).
#  - category: clear rows 3-6 leftovers.
#  - URL: swap the old "C:\" entry for a www.google.de link, clear
#    rows 3-6 leftovers and their hyperlinks.
#  - color: clear rows 3-6 leftovers.
#  - add a new "test" sheet with a tooltip-text value.
#  - mirror the selections/active tab recorded in the authored session.

$wb = $excel.ActiveWorkbook

$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsCategory   = $wb.Worksheets.Item("category")
$wsUrl        = $wb.Worksheets.Item("URL")
$wsColor      = $wb.Worksheets.Item("color")

# ---------------------------------------------------------------
# fieldnames
# ---------------------------------------------------------------
$wsFieldnames.Range("A1").Value = "Internet #1"
$wsFieldnames.Range("B1").Value = "Internet #2"
$wsFieldnames.Range("C1").Value = "Internet #3"
$wsFieldnames.Range("D1").Value = "Internet #4"

$wsFieldnames.Range("A2").Value = "wikipedia"
$wsFieldnames.Range("B2").Value = "Github"
$wsFieldnames.Range("C2").Value = "Expasy"
$wsFieldnames.Range("D2").Value = "duckduckgo"

$wsFieldnames.Range("A3:D6").ClearContents() | Out-Null

$wsFieldnames.Hyperlinks.Delete() | Out-Null
$wsFieldnames.Hyperlinks.Add($wsFieldnames.Range("C2"), "https://www.expasy.org/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Expasy") | Out-Null

# ---------------------------------------------------------------
# category
# ---------------------------------------------------------------
$wsCategory.Range("A3:D6").ClearContents() | Out-Null

# ---------------------------------------------------------------
# URL
# ---------------------------------------------------------------
$wsUrl.Range("A2").Value = "www.google.de"

$wsUrl.Range("A3:D6").ClearContents() | Out-Null

$wsUrl.Hyperlinks.Delete() | Out-Null
$wsUrl.Hyperlinks.Add($wsUrl.Range("A2"), "https://www.google.de/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "www.google.de") | Out-Null
$wsUrl.Hyperlinks.Add($wsUrl.Range("B2"), "http://www.github.com/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "www.github.com") | Out-Null
$wsUrl.Hyperlinks.Add($wsUrl.Range("D2"), "https://192.168.178.1/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://192.168.178.1/") | Out-Null

# ---------------------------------------------------------------
# color
# ---------------------------------------------------------------
$wsColor.Range("A3:D6").ClearContents() | Out-Null

# ---------------------------------------------------------------
# new "test" sheet, appended after "color"
# ---------------------------------------------------------------
$wsTest = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsColor)
$wsTest.Name = "test"
$wsTest.Range("A2").Value = "asdfsafd"

# ---------------------------------------------------------------
# replay the recorded selections on every sheet, then land on "test"
# ---------------------------------------------------------------
$wsFieldnames.Activate() | Out-Null
$wsFieldnames.Range("A2").Select() | Out-Null

$wsCategory.Activate() | Out-Null
$wsCategory.Range("D60").Select() | Out-Null

$wsUrl.Activate() | Out-Null
$wsUrl.Range("A3").Select() | Out-Null

$wsColor.Activate() | Out-Null
$wsColor.Range("G47").Select() | Out-Null

$wsTest.Activate() | Out-Null
$wsTest.Range("C8").Select() | Out-Null

Write-Host "edit applied"
